# Update loading_percent values for rows 2-25 (case with 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to: C, D, E, F, G, H, J, L, M, N, O
# (columns B, I, K are left unchanged at 0)
$colIndices = @(3, 4, 5, 6, 7, 8, 10, 12, 13, 14, 15)

# New values per row (rows 2 through 25), one array per row matching $colIndices order
$rowsData = @(
    @(12.78096850132919, 8.567211726214827, 13.37925413994992, 35.49796562664235, 39.87842566722569, 16.74819572065282, 10.18790123844742, 9.379824911065869, 25.20222021013713, 17.60445428078719, 27.06915098515913),
    @(12.80776577770738, 8.591611068104703, 13.43751292802409, 35.44780097638129, 39.62326520053956, 16.76234810793886, 10.22417599341521, 9.399967631910847, 24.5889098833892, 17.32554009187516, 27.03533672270149),
    @(12.82711724026596, 8.607583242300826, 13.47518103596884, 35.42632185143687, 39.47850878320243, 16.77442716687887, 10.24754616311255, 9.412970364607618, 24.20456403080352, 17.15397310621729, 27.02133908307387),
    @(12.83572980011269, 8.614341476803936, 13.49100891864269, 35.41991745426282, 39.42256737530149, 16.78020097161333, 10.25734634046113, 9.41842921351752, 24.0461784225695, 17.08406341102513, 27.01733889876572),
    @(12.83720373718737, 8.615478750575233, 13.4936660099304, 35.41899594033056, 39.41346385480075, 16.78121110585948, 10.2589903820007, 9.419345335855425, 24.01977805729469, 17.07245774944215, 27.01677764600116),
    @(12.82723045272002, 8.607673375849325, 13.47539256056479, 35.4262259662909, 39.47774193179418, 16.77450158811067, 10.24767721067425, 9.413043335603614, 24.20243485730677, 17.15303014674073, 27.02127823306567),
    @(12.78960560360476, 8.57541913015455, 13.39894854272207, 35.47873729564022, 39.78800692998449, 16.75237158253144, 10.20018145534939, 9.38663859373853, 24.9924757113226, 17.50839942411793, 27.05609007911984),
    @(12.73889943840674, 8.520018174125672, 13.26405602631551, 35.6553945270191, 40.48838828332968, 16.73589399196969, 10.11571612684069, 9.339876907898237, 26.47220084127454, 18.19919129865553, 27.17783668396317),
    @(12.71581743220674, 8.484082756965016, 13.17404985289227, 35.82958846153839, 41.05520802233175, 16.74021789712021, 10.05890010955943, 9.308550946240665, 27.50723448302654, 18.69832008517655, 27.29954090218512),
    @(12.708412496028, 8.468766409897341, 13.13506836793117, 35.918321194141, 41.32350847522618, 16.74575013197978, 10.0341809659551, 9.294951625992832, 27.96496492250514, 18.92262260590642, 27.36180763097835),
    @(12.70605466758258, 8.463114489080574, 13.12058856654357, 35.95327023724846, 41.42652804679154, 16.74835687964431, 10.0249817593088, 9.289895067015674, 28.13627482650153, 19.00708982982884, 27.38636786328725),
    @(12.70654260420564, 8.464325149840315, 13.12369453418365, 35.94568367939696, 41.40427913463125, 16.74777272459625, 10.02695580482976, 9.290979948439846, 28.09947219282765, 18.98892042302869, 27.38103493720948),
    @(12.70820956639012, 8.46829845693585, 13.13387146483732, 35.92116956363304, 41.33195598358319, 16.74595433947927, 10.03342091063588, 9.294533754126833, 27.97909997314784, 18.92958163897295, 27.36380862251241),
    @(12.70928877626073, 8.470751496919217, 13.14014178911643, 35.90632895656087, 41.28783838154828, 16.74490714704722, 10.03740197468259, 9.296722689190496, 27.90510122851232, 18.89317137229921, 27.35338442108755),
    @(12.71636375080429, 8.485104443467806, 13.17663682858402, 35.8239791542677, 41.03787725083129, 16.73992802846767, 10.06053818289589, 9.309452759063683, 27.47704488073961, 18.68359946037656, 27.29560958204465),
    @(12.72149767602038, 8.494173402010782, 13.19952763091952, 35.7758797289674, 40.88715307282288, 16.73778614335898, 10.07501963505267, 9.31742869507878, 27.21098656064389, 18.55427392101352, 27.26192688621381),
    @(12.72474185078616, 8.499486666878665, 13.21287859540882, 35.74910869995094, 40.80144941702569, 16.73688984188563, 10.08345507522824, 9.322077543415928, 27.05673017124096, 18.47963439202577, 27.24320399794593),
    @(12.72589026024938, 8.501302319575952, 13.21743075754263, 35.74019860002127, 40.77260395062978, 16.73664403935921, 10.08632941099161, 9.323662104531452, 27.0042953014449, 18.45432130233248, 27.23697680058892),
    @(12.72092100564586, 8.493197954023399, 13.19707174628073, 35.78090753458586, 40.90309616134706, 16.73797941688978, 10.07346708419811, 9.316573301613353, 27.23943691159885, 18.56806778315498, 27.26544521997883),
    @(12.70770781953511, 8.467127383833557, 13.13087461466509, 35.9283335175922, 41.35316118279938, 16.74647456156463, 10.03151757706037, 9.293487388827932, 28.01451214220273, 18.94702425097901, 27.36884187820995),
    @(12.70167377301273, 8.450951609659493, 13.0892523160745, 36.03253211943976, 41.65554117716869, 16.75500919919954, 10.00504160516801, 9.278942513849051, 28.50922346367613, 19.1919148809861, 27.44213030097741),
    @(12.70465588529902, 8.459506034501061, 13.11131693029172, 35.9762074660286, 41.49342977150006, 16.75018157606611, 10.01908648367447, 9.286655827010502, 28.24631286634881, 19.06149040431164, 27.40249632164651),
    @(12.72118080695631, 8.493638644427044, 13.19818145753023, 35.77863171502661, 40.89588532646371, 16.73789099399726, 10.0741686501676, 9.316959827544686, 27.22657853382017, 18.56183247108967, 27.2638525813622),
    @(12.75013484831817, 8.534167098857372, 13.29894628550205, 35.59976038902224, 40.28944648992082, 16.7374664331337, 10.13764254632204, 9.351993010428934, 26.0803920586682, 18.01345874597482, 27.13920544786342)
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $rowsData[$i]
    for ($j = 0; $j -lt $colIndices.Count; $j++) {
        $colNum = $colIndices[$j]
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

Write-Host "Updated loading_percent values for rows 2-25"
